$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9.663483666666666
$ws.Range("H2").Value = 28.990451
$ws.Range("I2").Value = 0.4172798466714015
$ws.Range("J2").Value = 0.4172798466714016
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 18.444833
$ws.Range("N2").Value = 55.33449900000001
$ws.Range("O2").Value = 0.529296397589589
$ws.Range("P2").Value = 0.5292963975895891
$ws.Range("Q2").Value = 178.2413424298944
$ws.Range("R2").Value = 1604.172081869049
$ws.Range("S2").Value = 0.2208647196299089
$ws.Range("T2").Value = 0.2208647196299089

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 9.663483666666666
$ws.Range("H3").Value = 28.990451
$ws.Range("I3").Value = 0.4172798466714015
$ws.Range("J3").Value = 0.4172798466714016
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 8.028768
$ws.Range("N3").Value = 24.086304
$ws.Range("O3").Value = 0.2303950368909585
$ws.Range("P3").Value = 0.2303950368909585
$ws.Range("Q3").Value = 77.58586843145599
$ws.Range("R3").Value = 698.2728158831039
$ws.Range("S3").Value = 0.09613920566771106
$ws.Range("T3").Value = 0.09613920566771107

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 9.663483666666666
$ws.Range("H4").Value = 28.990451
$ws.Range("I4").Value = 0.4172798466714015
$ws.Range("J4").Value = 0.4172798466714016
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 8.374233
$ws.Range("N4").Value = 25.122699
$ws.Range("O4").Value = 0.2403085655194523
$ws.Range("P4").Value = 0.2403085655194524
$ws.Range("Q4").Value = 80.92426381636099
$ws.Range("R4").Value = 728.318374347249
$ws.Range("S4").Value = 0.1002759213737815
$ws.Range("T4").Value = 0.1002759213737815

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.297564333333334
$ws.Range("H5").Value = 15.892693
$ws.Range("I5").Value = 0.2287546509102482
$ws.Range("J5").Value = 0.2287546509102482
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 18.444833
$ws.Range("N5").Value = 55.33449900000001
$ws.Range("O5").Value = 0.529296397589589
$ws.Range("P5").Value = 0.5292963975895891
$ws.Range("Q5").Value = 97.71268943508969
$ws.Range("R5").Value = 879.4142049158072
$ws.Range("S5").Value = 0.1210790126586584
$ws.Range("T5").Value = 0.1210790126586584

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.297564333333334
$ws.Range("H6").Value = 15.892693
$ws.Range("I6").Value = 0.2287546509102482
$ws.Range("J6").Value = 0.2287546509102482
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 8.028768
$ws.Range("N6").Value = 24.086304
$ws.Range("O6").Value = 0.2303950368909585
$ws.Range("P6").Value = 0.2303950368909585
$ws.Range("Q6").Value = 42.532914997408
$ws.Range("R6").Value = 382.796234976672
$ws.Range("S6").Value = 0.05270393623544498
$ws.Range("T6").Value = 0.05270393623544498

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.297564333333334
$ws.Range("H7").Value = 15.892693
$ws.Range("I7").Value = 0.2287546509102482
$ws.Range("J7").Value = 0.2287546509102482
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 8.374233
$ws.Range("N7").Value = 25.122699
$ws.Range("O7").Value = 0.2403085655194523
$ws.Range("P7").Value = 0.2403085655194524
$ws.Range("Q7").Value = 44.363038059823
$ws.Range("R7").Value = 399.267342538407
$ws.Range("S7").Value = 0.05497170201614484
$ws.Range("T7").Value = 0.05497170201614485

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 8.197232333333334
$ws.Range("H8").Value = 24.591697
$ws.Range("I8").Value = 0.3539655024183503
$ws.Range("J8").Value = 0.3539655024183503
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 18.444833
$ws.Range("N8").Value = 55.33449900000001
$ws.Range("O8").Value = 0.529296397589589
$ws.Range("P8").Value = 0.5292963975895891
$ws.Range("Q8").Value = 151.1965814505337
$ws.Range("R8").Value = 1360.769233054803
$ws.Range("S8").Value = 0.1873526653010218
$ws.Range("T8").Value = 0.1873526653010218

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 8.197232333333334
$ws.Range("H9").Value = 24.591697
$ws.Range("I9").Value = 0.3539655024183503
$ws.Range("J9").Value = 0.3539655024183503
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 8.028768
$ws.Range("N9").Value = 24.086304
$ws.Range("O9").Value = 0.2303950368909585
$ws.Range("P9").Value = 0.2303950368909585
$ws.Range("Q9").Value = 65.813676646432
$ws.Range("R9").Value = 592.323089817888
$ws.Range("S9").Value = 0.08155189498780248
$ws.Range("T9").Value = 0.08155189498780248

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 8.197232333333334
$ws.Range("H10").Value = 24.591697
$ws.Range("I10").Value = 0.3539655024183503
$ws.Range("J10").Value = 0.3539655024183503
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 8.374233
$ws.Range("N10").Value = 25.122699
$ws.Range("O10").Value = 0.2403085655194523
$ws.Range("P10").Value = 0.2403085655194524
$ws.Range("Q10").Value = 68.645533514467
$ws.Range("R10").Value = 617.809801630203
$ws.Range("S10").Value = 0.085060942129526
$ws.Range("T10").Value = 0.085060942129526

